$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.443.38"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "1.908.58"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4668"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4077"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08021"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "1.874.30"
$ws.Range("E13").Value = "  -1.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.942"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.125"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06594"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "29.450.91"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.537"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("D26").Value = "2.101.82"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.707"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("E32").Value = "  +9.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09472"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.576"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.387"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02255"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.376"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.173"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5869"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.304"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07749"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.379"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5539"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.923"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2936"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.21%  "
